# Daily update of Burndown Chart
# - Fill in today's "Actual" progress values for 7/27 - 7/30 (rows 5-8, column C)
#   to match the "Planned" values already present in column B.
# - Leave the active selection on E3, matching where the author clicked next.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C5").Value = 27
$ws.Range("C6").Value = 27
$ws.Range("C7").Value = 27
$ws.Range("C8").Value = 27

$ws.Range("E3").Select()
